$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) stays the same.

# Data rows 2-9: Sending cluster / Ligand symbol / Receptor symbol / Target cluster
# plus numeric columns E..T.
$data = @(
    @("FAPs","Lgi3","Adam23","ECs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,0.1448543333333333,0.434563,0.007514552731137001,0.007514552731137,0.2224630843576667,2.002167759219,0.00672939638314816,0.006729396383148159),
    @("FAPs","Lgi3","Adam23","FAPs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,14.494489,43.483467,0.7519250504625467,0.7519250504625465,22.26021586601901,200.341942794171,0.6733603310832776,0.6733603310832775),
    @("FAPs","Lgi3","Adam23","M2",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,2,0.6666666666666666,0.091182,0.273546,0.004730213666123443,0.004730213666123443,0.140034671322,1.260312041898,0.004235978357625124,0.004235978357625124),
    @("FAPs","Lgi3","Adam23","sCs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,4.545982333333334,13.637947,0.235830183140193,0.2358301831401929,6.981587834045668,62.83429050641101,0.2111895196217034,0.2111895196217034),
    @("sCs","Lgi3","Adam23","ECs",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,0.1448543333333333,0.434563,0.007514552731137001,0.007514552731137,0.025956013427,0.233604120843,0.0007851563479888403,0.0007851563479888402),
    @("sCs","Lgi3","Adam23","FAPs",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,14.494489,43.483467,0.7519250504625467,0.7519250504625465,2.597224000443,23.375016003987,0.07856471937926895,0.07856471937926894),
    @("sCs","Lgi3","Adam23","M2",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,2,0.6666666666666666,0.091182,0.273546,0.004730213666123443,0.004730213666123443,0.016338629034,0.147047661306,0.0004942353084983197,0.0004942353084983197),
    @("sCs","Lgi3","Adam23","sCs",1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,4.545982333333334,13.637947,0.235830183140193,0.2358301831401929,0.814580936363,7.331228427267,0.02464066351848952,0.02464066351848952)
)

$rowIdx = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $val
        $colIdx++
    }
    $rowIdx++
}
